$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Layouts we need from the slide master ("Office Theme"):
#   2 = "Title and Content"  (title + single body placeholder idx=1)
#   5 = "Comparison"         (title + body idx1 + half idx2 + body-quarter idx3 + quarter idx4)
# ---------------------------------------------------------------------------
$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)

# ---------------------------------------------------------------------------
# 1) The original slide 2 used the "Comparison" layout but was completely
#    empty. That slide (with its five placeholders) is being repurposed as
#    the new slide 5 ("Beginning Development") with text filled in, so grab
#    a duplicate of it *before* we touch slide 2's own content.
# ---------------------------------------------------------------------------
$origSlide2 = $p.Slides.Item(2)
$dupRange = $origSlide2.Duplicate()
$newSlide5 = $dupRange.Item(1)
$newSlide5.MoveTo($p.Slides.Count)

# Fill in the text for the new slide 5 ("Beginning Development").
$newSlide5.Shapes.Item(1).TextFrame.TextRange.Text = "Beginning Development"
$newSlide5.Shapes.Item(2).TextFrame.TextRange.Text = "Flowchart"
$newSlide5.Shapes.Item(4).TextFrame.TextRange.Text = "wireframe"

# ---------------------------------------------------------------------------
# 2) Turn the original slide 2 into a simple "Title and Content" slide
#    ("Value Added Proposal"). Easiest reliable way: delete it and insert a
#    fresh slide with the right layout in the same position.
# ---------------------------------------------------------------------------
$origSlide2.Delete()

$slide2 = $p.Slides.AddSlide(2, $titleContentLayout)
$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "Value Added Proposal"

# ---------------------------------------------------------------------------
# 3) Two brand-new "Title and Content" slides inserted right after it.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.AddSlide(3, $titleContentLayout)
$slide3.Shapes.Item(1).TextFrame.TextRange.Text = "Group Members"

$slide4 = $p.Slides.AddSlide(4, $titleContentLayout)
$slide4.Shapes.Item(1).TextFrame.TextRange.Text = "Tools Used for Development"

Write-Output ("Final slide count: " + $p.Slides.Count.ToString())
for ($i = 1; $i -le $p.Slides.Count; $i++) {
  $sl = $p.Slides.Item($i)
  $title = ""
  if ($sl.Shapes.HasTitle) {
    $title = $sl.Shapes.Title.TextFrame.TextRange.Text
  }
  Write-Output ("Slide " + $i.ToString() + " shapes=" + $sl.Shapes.Count.ToString() + " title='" + $title + "'")
}
